$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark from the end of the document to the
#        middle of the "Furthermore, this approach doesn't work..." run,
#        splitting that run into two pieces: "...doesn't work und" | "er the
#        current framework...".
$marker = ". Furthermore, this approach doesn" + [char]0x2019 + "t work und"
$full = $d.Content.Text
$splitPos = $full.IndexOf($marker) + $marker.Length

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null

# --- 2. Add a new, empty paragraph at the very end of the document (after
#        the "...top 3 interactors." paragraph, before the sectPr).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
